$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1. Column A width -> narrow (renders as width=1 in OOXML units) ---
$ws.Columns("A").ColumnWidth = 0.14

# --- 2. Header row 2 (A2:W2): restyle from the "L0..Activation" sub-header look
#        to the bold/grey style already used by row 1 (style index 10) ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:W2").PasteSpecial(-4122) | Out-Null

# --- 3. Update GEMM (U) totals for the stride-2 conv rows (back-prop of a
#        strided convolution costs 4x the forward GEMM) ---
$uRows = @"
3,472055808
46,462422016
52,411041792
90,462422016
96,411041792
154,462422016
160,411041792
"@
foreach ($line in ($uRows -split "`r?`n")) {
    if ($line.Trim() -eq "") { continue }
    $parts = $line -split ","
    $r = [int]$parts[0]
    $v = [double]$parts[1]
    $ws.Cells.Item($r, 21).Value = $v
}

# --- 4. New column X: "BackProp" header + per-row totals ---
# First extend the merged "Operation Summary" header box from U1:W1 to U1:X1.
# W1 stops being the closing cell of the box (loses its right border) and
# X1 becomes the new closing cell (gains it), so re-use the existing
# "middle"/"end of merge" styles already present at V1 / W1.
$ws.Range("U1:W1").UnMerge() | Out-Null
$ws.Range("V1").Copy() | Out-Null
$ws.Range("W1").PasteSpecial(-4122) | Out-Null
$ws.Range("W1").Copy() | Out-Null
$ws.Range("X1").PasteSpecial(-4122) | Out-Null
$ws.Range("U1:X1").Merge() | Out-Null

$ws.Cells.Item(2, 24).Value = "BackProp"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("X2").PasteSpecial(-4122) | Out-Null

$xRows = @"
3,472055808
4,1605632
5,802816
6,2517630976
7,
8,
9,12845056
10,401408
11,200704
12,115605504
13,401408
14,200704
15,51380224
16,1605632
17,
18,51380224
19,1605632
20,802816
21,
22,51380224
23,401408
24,200704
25,115605504
26,401408
27,200704
28,51380224
29,1605632
30,802816
31,
32,51380224
33,401408
34,200704
35,115605504
36,401408
37,200704
38,51380224
39,1605632
40,802816
41,
42,
43,102760448
44,802816
45,401408
46,462422016
47,200704
48,100352
49,51380224
50,802816
51,
52,411041792
53,802816
54,401408
55,
56,51380224
57,200704
58,100352
59,115605504
60,200704
61,100352
62,51380224
63,802816
64,401408
65,
66,51380224
67,200704
68,100352
69,115605504
70,200704
71,100352
72,51380224
73,802816
74,401408
75,
76,51380224
77,200704
78,100352
79,115605504
80,200704
81,100352
82,51380224
83,802816
84,401408
85,
86,
87,102760448
88,401408
89,200704
90,462422016
91,100352
92,50176
93,51380224
94,401408
95,
96,411041792
97,401408
98,200704
99,
100,51380224
101,100352
102,50176
103,115605504
104,100352
105,50176
106,51380224
107,401408
108,200704
109,
110,51380224
111,100352
112,50176
113,115605504
114,100352
115,50176
116,51380224
117,401408
118,200704
119,
120,51380224
121,100352
122,50176
123,115605504
124,100352
125,50176
126,51380224
127,401408
128,200704
129,
130,51380224
131,100352
132,50176
133,115605504
134,100352
135,50176
136,51380224
137,401408
138,200704
139,
140,51380224
141,100352
142,50176
143,115605504
144,100352
145,50176
146,51380224
147,401408
148,200704
149,
150,
151,102760448
152,200704
153,100352
154,462422016
155,50176
156,25088
157,51380224
158,200704
159,
160,411041792
161,200704
162,100352
163,
164,51380224
165,50176
166,25088
167,115605504
168,50176
169,25088
170,51380224
171,200704
172,100352
173,
174,51380224
175,50176
176,25088
177,115605504
178,50176
179,25088
180,51380224
181,200704
182,100352
183,
184,2048000
"@
foreach ($line in ($xRows -split "`r?`n")) {
    if ($line.Trim() -eq "") { continue }
    $parts = $line -split ","
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 24).Value2 = $null
    if ($parts.Length -gt 1 -and $parts[1].Trim() -ne "") {
        $v = [double]$parts[1]
        $ws.Cells.Item($r, 24).Value = $v
    }
    $ws.Range("W" + $r).Copy() | Out-Null
    $ws.Range("X" + $r).PasteSpecial(-4122) | Out-Null
}

Write-Output "stage2 done"
